$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 - Fortuna Sittard
$ws.Range("C12").Value = 24
$ws.Range("D12").Value = 8
$ws.Range("G12").Value = 36
$ws.Range("H12").Value = 43
$ws.Range("I12").Value = -7
$ws.Range("J12").Value = 29
$ws.Range("K12").Value = 1.21
$ws.Range("L12").Value = "W L D L W"
$ws.Range("M12").Value = 10558
$ws.Range("N12").Value = "Kaj Sierhuis - 9"

# Row 13 - Excelsior
$ws.Range("C13").Value = 24
$ws.Range("F13").Value = 12
$ws.Range("G13").Value = 26
$ws.Range("H13").Value = 41
$ws.Range("I13").Value = -15
$ws.Range("K13").Value = 1.08
$ws.Range("L13").Value = "D D W L L"
$ws.Range("N13").Value = "Noah Naujoks - 7"
